$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 text: remove "RME" from the steel description
$ws.Range("B2").Value = "22% CR+PC/LFM+CDL/H:1`n5% CR+PC/LFM+CDH/H:2`n5% CR+PC/LFM+CDM/H:2`n30% S/LFM+CDL/H:1`n5% S/LFM+CDH/H:1`n5% S/LFM+CDM/H:1`n6% S+SL/LFM+CDL/H:1`n3% S+SL/LFM+CDH/H:1`n5% S/LFBR+CDH/H:1`n5% CR/LFM+CDL/H:2`n1% CR/LFM+CDM/H:2`n4% MCF/LWAL+CDL/H:1`n4% MUR/LWAL+CDN/H:1"

# Apply wrap text formatting to B2 and set row height
$ws.Range("B2").WrapText = $true
$ws.Rows(2).RowHeight = 409.6

# Update the selection to match the saved view state
$ws.Range("E2:E8").Select() | Out-Null
